# feat: add 2022-Q1 data
#
# 1. Create a new "2022-Q1" sheet (positioned right after "2021-Q4", before
#    "总计") by duplicating the "2021-Q4" sheet - this keeps the existing
#    header row / cell styling byte-identical - then overwrite the data row
#    with the 2022-Q1 numbers.
# 2. Insert a new top data row into the "总计" sheet for "2022-Q1" and shift
#    the existing quarters down, renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# Header row (B1:H1) is identical to "2021-Q4"'s header, so only the data
# row (row 2) needs updating. B2/C2 (fund code/name) are unchanged too.
$new.Range("D2:G2").NumberFormat = "@"
$new.Range("D2").Value = "3.37"
$new.Range("E2").Value = "95.08"
$new.Range("F2").Value = "2.73"
$new.Range("G2").Value = "0.0920"
$new.Range("D2:G2").ClearFormats()
$new.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 2) Update "总计" sheet: insert a new row for 2022-Q1 above the existing
#    data and renumber the index column.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Restore the index-column style (lost a border component during the row
# insert) by pulling formats back from the row right below, which still
# carries the original style.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.09

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
